$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Range("F2").Value = 38
$ws.Range("F4").Value = 170
$ws.Range("F5").Value = 490
$ws.Range("F6").Value = 844
$ws.Range("F7").Value = 265
$ws.Range("F8").Value = 1248
$ws.Range("F9").Value = 365
$ws.Range("F11").Value = 894
$ws.Range("F12").Value = 713
$ws.Range("F15").Value = 145
$ws.Range("F18").Value = 2980
$ws.Range("F19").Value = 2635
$ws.Range("F23").Value = 317
$ws.Range("F26").Value = 5355
$ws.Range("F30").Value = 61
$ws.Range("F32").Value = 1117
$ws.Range("F33").Value = 72
$ws.Range("F35").Value = 299

$ws = $wb.Worksheets.Item("演出")
$ws.Range("F4").Value = 1147
$ws.Range("F7").Value = 235
$ws.Range("F9").Value = 331
$ws.Range("F14").Value = 617
$ws.Range("F17").Value = 990
$ws.Range("F24").Value = 323
$ws.Range("F25").Value = 283
$ws.Range("F26").Value = 3978
$ws.Range("F30").Value = 201

$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F2").Value = 1787
$ws.Range("F5").Value = 2498
$ws.Range("F6").Value = 1062
$ws.Range("F9").Value = 1353
$ws.Range("F10").Value = 371

$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F2").Value = 1787
$ws.Range("F5").Value = 2498
$ws.Range("F7").Value = 1062
$ws.Range("F8").Value = 1353
$ws.Range("F9").Value = 371
$ws.Range("F11").Value = 170
$ws.Range("F12").Value = 490
$ws.Range("F13").Value = 844
$ws.Range("F14").Value = 265
$ws.Range("F15").Value = 1248
$ws.Range("F16").Value = 365
$ws.Range("F17").Value = 894
$ws.Range("F18").Value = 713
$ws.Range("F19").Value = 1147
$ws.Range("F20").Value = 1147
$ws.Range("F23").Value = 145
$ws.Range("F25").Value = 2980
$ws.Range("F26").Value = 2635
$ws.Range("F28").Value = 317
$ws.Range("F32").Value = 5355
$ws.Range("F35").Value = 617
$ws.Range("F36").Value = 617
$ws.Range("F38").Value = 61
$ws.Range("F44").Value = 323
$ws.Range("F45").Value = 323
$ws.Range("F46").Value = 1117
$ws.Range("F47").Value = 201
$ws.Range("F51").Value = 299

